# This workbook's data rows (2-7) were reordered (their full content, not
# just individual fields, moved to different row positions). Reconstruct
# the new order by first snapshotting each row's full contents, then
# writing them back to their new positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns Y and AA hold date-like text (e.g. "2023-08-11") stored as plain
# text in the original file. Writing such strings back via Value2 would
# normally be auto-converted by Excel into a date serial number, so force
# a text number format on those columns first to preserve them as text.
$ws.Range("Y2:Y7").NumberFormat = "@"
$ws.Range("AA2:AA7").NumberFormat = "@"

# Snapshot full row contents (columns A through AY) for all data rows.
$rowsData = @{}
for ($r = 2; $r -le 7; $r++) {
    $rng = $ws.Range("A" + $r + ":AY" + $r)
    $rowsData[$r] = $rng.Value2
}

# Mapping of new row number -> row number whose original content should
# now occupy it.
$mapping = @{2=3; 3=7; 4=2; 5=6; 6=5; 7=4}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $destRng = $ws.Range("A" + $newRow + ":AY" + $newRow)
    $destRng.Value2 = $rowsData[$oldRow]
}
